$wb = $excel.ActiveWorkbook

# Existing sheets: "总计" (summary) and "2022-Q3"
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item(2)

# 1) Insert the new "2022-Q4" worksheet between "总计" and "2022-Q3" by
#    duplicating the summary sheet (this keeps the same sheetPr/pageSetup
#    structure and cell style "2") and then replacing its contents.
$wsTotal.Copy($wsQ3, $null) | Out-Null
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"
$wsQ4.Cells.ClearContents()

# Extend the header/row style (style index 2, already used on $wsTotal) to
# the extra columns E:H and to row 3 column A.
$wsTotal.Range("B1").Copy() | Out-Null
$wsQ4.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$wsTotal.Range("A2").Copy() | Out-Null
$wsQ4.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2) Populate the new "2022-Q4" sheet with fund holding data
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2:G2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "005189"
$wsQ4.Range("C2").Value = "海富通量化前锋股票A"
$wsQ4.Range("D2").Value = "0.54"
$wsQ4.Range("E2").Value = "88.13"
$wsQ4.Range("F2").Value = "1.00"
$wsQ4.Range("G2").Value = "0.0054"
$wsQ4.Range("B2:G2").Style = "Normal"
$wsQ4.Range("H2").Value = 7

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3:G3").NumberFormat = "@"
$wsQ4.Range("B3").Value = "005188"
$wsQ4.Range("C3").Value = "海富通量化前锋股票C"
$wsQ4.Range("D3").Value = "0.03"
$wsQ4.Range("E3").Value = "88.13"
$wsQ4.Range("F3").Value = "1.00"
$wsQ4.Range("G3").Value = "0.0003"
$wsQ4.Range("B3:G3").Style = "Normal"
$wsQ4.Range("H3").Value = 7

# 3) Update the "总计" (summary) sheet: shift the old 2022-Q3 row down to row 3,
#    and write the new 2022-Q4 totals into row 2.
$wsTotal.Range("A2").Copy() | Out-Null
$wsTotal.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.06

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01
